$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.792.26'
$ws.Range('E2').Value = '  -1.73%  '
$ws.Range('D3').Value = '1.869.07'
$ws.Range('E3').Value = '  -2.00%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '300.06'
$ws.Range('E5').Value = '  -2.51%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').Value = '0.5333'
$ws.Range('E7').Value = '  +1.29%  '
$ws.Range('D8').Value = '0.3738'
$ws.Range('E8').Value = '  -2.17%  '
$ws.Range('D9').Value = '0.07142'
$ws.Range('E9').Value = '  -2.10%  '
$ws.Range('D10').Value = '21.55'
$ws.Range('E10').Value = '  -2.51%  '
$ws.Range('D11').Value = '0.8869'
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').Value = '0.08140'
$ws.Range('E12').Value = '  -0.58%  '
$ws.Range('D13').Value = '1.898.70'
$ws.Range('E13').Value = '  +37.35%  '
$ws.Range('E14').Value = '  -3.68%  '
$ws.Range('D15').Value = '5.288'
$ws.Range('E15').Value = '  -1.36%  '
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('D18').Value = '0.000008483'
$ws.Range('E18').Value = '  -1.84%  '
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').Value = '26.822.70'
$ws.Range('E20').Value = '  -1.79%  '
$ws.Range('D21').Value = '4.968'
$ws.Range('E22').Value = '  -1.72%  '
$ws.Range('E23').Value = '  -2.27%  '
$ws.Range('D24').Value = '2.285'
$ws.Range('E24').Value = '  -0.75%  '
$ws.Range('D25').Value = '145.97'
$ws.Range('D26').Value = '1.742'
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  -1.46%  '
$ws.Range('D28').Value = '113.68'
$ws.Range('E28').Value = '  -2.51%  '
$ws.Range('D29').Value = '4.694'
$ws.Range('E29').Value = '  -2.94%  '
$ws.Range('D30').Value = '4.623'
$ws.Range('E30').Value = '  -4.37%  '
$ws.Range('D31').Value = '0.09102'
$ws.Range('E31').Value = '  -1.73%  '
$ws.Range('D32').Value = '0.8085'
$ws.Range('E32').Value = '  -2.85%  '
$ws.Range('D33').Value = '0.05018'
$ws.Range('E33').Value = '  -1.27%  '
$ws.Range('D34').Value = '1.170'
$ws.Range('E34').Value = '  -5.04%  '
$ws.Range('D35').Value = '2.946'
$ws.Range('E35').Value = '  -1.84%  '
$ws.Range('D36').Value = '0.6119'
$ws.Range('E36').Value = '  +4.95%  '
$ws.Range('D37').Value = '2.693'
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('D38').Value = '3.180'
$ws.Range('E38').Value = '  -5.38%  '
$ws.Range('D39').Value = '0.01943'
$ws.Range('E39').Value = '  -3.17%  '
$ws.Range('E40').Value = '  -1.52%  '
$ws.Range('D41').Value = '0.5271'
$ws.Range('E41').Value = '  +6.94%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = '8.776'
$ws.Range('E42').Value = '  -5.94%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '6.477'
$ws.Range('E43').Value = '  -1.15%  '
$ws.Range('D44').Value = '116.22'
$ws.Range('E44').Value = '  -0.62%  '
$ws.Range('E45').Value = '  -2.55%  '
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').Value = '1.646'
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('D48').Value = '9.943'
$ws.Range('E48').Value = '  -2.13%  '
$ws.Range('D49').Value = '37.25'
$ws.Range('E49').Value = '  -4.14%  '
$ws.Range('D50').Value = '0.06066'
$ws.Range('E50').Value = '  -2.26%  '
$ws.Range('D51').Value = '62.17'
$ws.Range('E51').Value = '  -2.98%  '
